# Auto-generated edit script applying the cryptos price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = "D4","D5","D6","D7","D8","D9","D11","D12","D13","D14","D15","D17","D18","D19","D20","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51"
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "30.480.80"
$ws.Range("D3").Value = "1.900.87"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "238.91"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.4898"
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("D8").Value = "0.2921"
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("D9").Value = "0.06676"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "1.929.61"
$ws.Range("E10").Value = "  +2.88%  "
$ws.Range("D11").Value = "16.99"
$ws.Range("E11").Value = "  +2.23%  "
$ws.Range("D12").Value = "0.07326"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D13").Value = "5.195"
$ws.Range("E13").Value = "  +3.86%  "
$ws.Range("D14").Value = "88.28"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "0.6681"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").Value = "30.454.98"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "13.44"
$ws.Range("E17").Value = "  +3.22%  "
$ws.Range("D18").Value = "0.000007864"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "5.456"
$ws.Range("E20").Value = "  +15.23%  "
$ws.Range("D21").Value = "2.152.69"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "196.96"
$ws.Range("E23").Value = "  -7.08%  "
$ws.Range("D24").Value = "6.155"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").Value = "9.470"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").Value = "163.27"
$ws.Range("E26").Value = "  +4.21%  "
$ws.Range("D27").Value = "18.40"
$ws.Range("E27").Value = "  -4.03%  "
$ws.Range("D28").Value = "1.940"
$ws.Range("E28").Value = "  +6.01%  "
$ws.Range("D29").Value = "1.472"
$ws.Range("E29").Value = "  +4.05%  "
$ws.Range("D30").Value = "4.332"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("D31").Value = "0.09170"
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("D32").Value = "4.136"
$ws.Range("E32").Value = "  +5.30%  "
$ws.Range("D33").Value = "0.05173"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").Value = "0.7398"
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").Value = "1.111"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("D36").Value = "2.728"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").Value = "0.01846"
$ws.Range("E37").Value = "  +1.74%  "
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").Value = "0.9234"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").Value = "2.067"
$ws.Range("E40").Value = "  +1.40%  "
$ws.Range("D41").Value = "0.4411"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "107.15"
$ws.Range("E42").Value = "  +2.51%  "
$ws.Range("D43").Value = "5.906"
$ws.Range("E43").Value = "  +2.91%  "
$ws.Range("D44").Value = "0.9950"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "69.08"
$ws.Range("E45").Value = "  +21.05%  "
$ws.Range("D46").Value = "0.1373"
$ws.Range("E46").Value = "  +3.47%  "
$ws.Range("D47").Value = "7.564"
$ws.Range("E47").Value = "  +3.07%  "
$ws.Range("D48").Value = "9.036"
$ws.Range("E48").Value = "  +5.13%  "
$ws.Range("D49").Value = "34.97"
$ws.Range("E49").Value = "  +5.23%  "
$ws.Range("D50").Value = "0.05831"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").Value = "0.3931"
$ws.Range("E51").Value = "  -2.18%  "
